# Apply the "Changes in the Script" commit to the CreateOffers sheet:
#  - Insert a new "Start_Date" column (C) between Description and Expiry_Date
#  - Update the existing offer row's Start_Date/Expiry_Date values
#  - Add a new offer row (Iseek2 / Testing3 / ...)
#  - Add a new blank (but formatted) row below it
#  - Update the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CreateOffers")

# Insert a new column before the old "Expiry_Date" column (col C) for "Start_Date"
$ws.Columns.Item(3).Insert()

# Give the new column roughly the same look as its neighbours (not bestFit)
$ws.Columns.Item(3).ColumnWidth = 10.33

# Header for the new column
$ws.Range("C1").Value = "Start_Date"

# Existing data row: split the old Expiry_Date value into Start_Date / Expiry_Date
$ws.Range("C2").Value = "2/12/2019"
$ws.Range("D2").Value = "12/31/2019"

# Copy the formatting (fill/border/style) of row 2 down into the new rows 3 and 4
$ws.Range("A2:G2").Copy()
$ws.Range("A3:G4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new offer row
$ws.Range("A3").Value = "Iseek2"
$ws.Range("B3").Value = "Testing3"
$ws.Range("C3").Value = "2/12/2019"
$ws.Range("D3").Value = "12/31/2019"
$ws.Range("E3").Value = "25"
$ws.Range("F3").Value = "2500"
$ws.Range("G3").Value = "Default Testing"

# Row 4 stays blank (only formatting was copied above)

# Match the saved selection/active cell
$ws.Range("F7").Select() | Out-Null
